$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new quarterly sheet "2022-Q4" positioned right before "2022-Q3"
#    (copy the "2022-Q3" sheet so it inherits the same layout/formatting,
#    then overwrite its contents with the 2022-Q4 fund-holdings data).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The template sheet has 4 data rows; 2022-Q4 only has 3, drop the extra one.
$q4.Rows.Item(5).Delete()

# Columns B:G store text (fund codes, names and formatted numeric strings),
# force text formatting first so values such as "000586" keep leading zeros.
$q4.Range("B2:G4").NumberFormat = "@"

$q4.Range("B2").Value = "000586"
$q4.Range("C2").Value = "景顺长城中小创精选股票"
$q4.Range("D2").Value = "2.30"
$q4.Range("E2").Value = "89.60"
$q4.Range("F2").Value = "8.22"
$q4.Range("G2").Value = "0.1891"
$q4.Range("H2").Value = 3

$q4.Range("B3").Value = "260115"
$q4.Range("C3").Value = "景顺长城中小盘混合"
$q4.Range("D3").Value = "1.06"
$q4.Range("E3").Value = "91.71"
$q4.Range("F3").Value = "4.84"
$q4.Range("G3").Value = "0.0513"
$q4.Range("H3").Value = 6

$q4.Range("B4").Value = "002597"
$q4.Range("C4").Value = "兴业成长动力灵活配置混合"
$q4.Range("D4").Value = "1.68"
$q4.Range("E4").Value = "89.03"
$q4.Range("F4").Value = "2.31"
$q4.Range("G4").Value = "0.0388"
$q4.Range("H4").Value = 6

# Drop the helper text-format style again now that the values are entered as
# text, so the cells end up with the same "no explicit style" look as the
# rest of the data rows on this sheet.
$q4.Range("B2:G4").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) summary sheet: insert a new row for 2022-Q4
#    above the existing 2022-Q3 entry, shifting everything else down.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()

# Insert() leaves the new row carrying formatting bled in from the header
# row above; strip it so the new data row matches the unstyled look of the
# other data rows, then reapply just the index-column style (copied from the
# row below, which still carries the correct style for column A).
$tot.Range("A2:D2").ClearFormats()
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("C2").Value = 3
$tot.Range("D2").Value = 0.28

# Column A is a plain 0-based row index; renumber every row below the newly
# inserted one now that they have all shifted down by one position.
for ($r = 3; $r -le 10; $r++) {
    $tot.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------------
# 3. Keep the tab-selection on the last sheet ("2020-Q4"), matching the
#    original workbook state (inserting/copying sheets can move it).
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item("2020-Q4")
$last.Activate()
